$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing F1 (2), G1 (781228) and H1 ("np") values to make room
# for two new numeric columns. Use .Value2 for the reads (this runtime's
# .Value getter doesn't resolve cleanly as an rvalue) and write back to
# front so nothing is clobbered before it's copied onward.
$ws.Range("L1").Value = $ws.Range("H1").Value2
$ws.Range("J1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = $ws.Range("F1").Value2

# Newly inserted numeric cells.
$ws.Range("F1").Value = 1.65
$ws.Range("H1").Value = 9.56
$ws.Range("I1").Value = 1526

# New text cell K1 - force Text number format so the date-looking string
# is stored verbatim (as a shared string), not coerced to a date serial.
$ws.Range("K1").NumberFormat = "@"
$ws.Range("K1").Value = "02/03/2024"

# Column K width (explicit width from the diff).
$ws.Columns.Item(11).ColumnWidth = 10.7109375

# Selection now sits on K1.
$ws.Range("K1").Select()
